$wb = $excel.ActiveWorkbook

$doc = $wb.Worksheets.Item("Document")
$db = $wb.Worksheets.Item("Database")

# Set new data values
$db.Range("D5").Value = 1
$db.Range("D6").Value = 1

# Update selections on each sheet
$doc.Range("B27").Select()
$db.Activate()
$db.Range("G7").Select()
